$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 5 - 保險 (insurance): widen from B:E to B:K, adding the common
# trailing "category/date/legislator_name/legislator_id/source_file/index"
# columns (F:K) shared by the other sheets, plus a "species" column (E)
# and renaming/reordering the "company"/name columns (B,C).
# ---------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item(5)

# The "date" column holds an ISO-looking string ("2013-12-31"); format the
# cells as Text first so Excel keeps it as a literal string instead of
# auto-converting it to a date serial number.
$wsIns.Range("G2:G4").NumberFormat = "@"

# Header row (row 1)
$wsIns.Range("B1").Value = "company"
$wsIns.Range("C1").Value = "name"
$wsIns.Range("D1").Value = "owner"
$wsIns.Range("E1").Value = "property_category"
$wsIns.Range("F1").Value = "category"
$wsIns.Range("G1").Value = "date"
$wsIns.Range("H1").Value = "legislator_name"
$wsIns.Range("I1").Value = "legislator_id"
$wsIns.Range("J1").Value = "source_file"
$wsIns.Range("K1").Value = "index"

# Row 2 (index 101)
$wsIns.Range("B2").Value = "國泰人壽"
$wsIns.Range("C2").Value = "鍾愛一生313"
$wsIns.Range("D2").Value = "黃靜秋"
$wsIns.Range("E2").Value = "insurance"
$wsIns.Range("F2").Value = "normal"
$wsIns.Range("G2").Value = "2013-12-31"
$wsIns.Range("H2").Value = "羅明才"
$wsIns.Range("I2").Value = 879
$wsIns.Range("J2").Value = "tmped981"
$wsIns.Range("K2").Value = 101

# Row 3 (index 102)
$wsIns.Range("B3").Value = "保德信國際人壽"
$wsIns.Range("C3").Value = "教育終身壽險"
$wsIns.Range("D3").Value = "黃靜秋"
$wsIns.Range("E3").Value = "insurance"
$wsIns.Range("F3").Value = "normal"
$wsIns.Range("G3").Value = "2013-12-31"
$wsIns.Range("H3").Value = "羅明才"
$wsIns.Range("I3").Value = 879
$wsIns.Range("J3").Value = "tmped981"
$wsIns.Range("K3").Value = 102

# Row 4 (index 103)
$wsIns.Range("B4").Value = "保德信國際人壽"
$wsIns.Range("C4").Value = "教育终身壽險"
$wsIns.Range("D4").Value = "黃靜秋"
$wsIns.Range("E4").Value = "insurance"
$wsIns.Range("F4").Value = "normal"
$wsIns.Range("G4").Value = "2013-12-31"
$wsIns.Range("H4").Value = "羅明才"
$wsIns.Range("I4").Value = 879
$wsIns.Range("J4").Value = "tmped981"
$wsIns.Range("K4").Value = 103

# ---------------------------------------------------------------------
# Sheet 6 - 債務 (debt): widen from B:G to B:N with the same trailing
# columns, rename header columns to species/debtor, and drop the old
# investment row (id 115) entirely - that data moved elsewhere.
# ---------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item(6)

# Same text-vs-date guard for the "date" column.
$wsDebt.Range("J2").NumberFormat = "@"

# Header row (row 1)
$wsDebt.Range("B1").Value = "species"
$wsDebt.Range("C1").Value = "debtor"
$wsDebt.Range("D1").Value = "owner"
$wsDebt.Range("E1").Value = "total"
$wsDebt.Range("F1").Value = "register_date"
$wsDebt.Range("G1").Value = "register_reason"
$wsDebt.Range("H1").Value = "property_category"
$wsDebt.Range("I1").Value = "category"
$wsDebt.Range("J1").Value = "date"
$wsDebt.Range("K1").Value = "legislator_name"
$wsDebt.Range("L1").Value = "legislator_id"
$wsDebt.Range("M1").Value = "source_file"
$wsDebt.Range("N1").Value = "index"

# Row 2 (index 113)
$wsDebt.Range("B2").Value = "房屋貸款"
$wsDebt.Range("C2").Value = "黃靜秋"
$wsDebt.Range("D2").Value = "台灣土地銀行臺北市中正區館前路46號"
$wsDebt.Range("E2").Value = 1294341
$wsDebt.Range("F2").Value = "89年10月25日"
$wsDebt.Range("G2").Value = "抵押"
$wsDebt.Range("H2").Value = "debt"
$wsDebt.Range("I2").Value = "normal"
$wsDebt.Range("J2").Value = "2013-12-31"
$wsDebt.Range("K2").Value = "羅明才"
$wsDebt.Range("L2").Value = 879
$wsDebt.Range("M2").Value = "tmped981"
$wsDebt.Range("N2").Value = 113

# Remove the old row 3 (investment entry, id 115) - no longer present
$wsDebt.Rows.Item(3).Delete()
